# Reorder the "Recorded By" (column G) entries on the
# "Session Analysis Results" sheet: for any cell whose value is a
# comma-separated list of recorders that starts with "System"/"system",
# reverse the order of the list (so "System" ends up last instead of first).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "

        if ($parts.Count -gt 1 -and $parts[0].ToLower() -eq "system") {
            $reversed = $parts[($parts.Count - 1)..0]
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
